# Applies scheduled market-data refresh to the Aegis_Profits sheets.
# Updates computed price/profit columns (H-N) for the affected leve rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 561.7143
$ws.Range("I29").Value = 283
$ws.Range("K29").Value = 849
$ws.Range("M29").Value = -568
$ws.Range("H38").Value = 2688366.5
$ws.Range("I38").Value = 3225922.8
$ws.Range("J38").Value = 585
$ws.Range("K38").Value = 9677768.399999999
$ws.Range("L38").Value = 1755
$ws.Range("M38").Value = -9677396.399999999
$ws.Range("N38").Value = -2499
$ws.Range("H43").Value = 1139.8
$ws.Range("I43").Value = 1139.8
$ws.Range("K43").Value = 1139.8
$ws.Range("M43").Value = -1070.8
$ws.Range("H58").Value = 1264205.9
$ws.Range("I58").Value = 2525411.8
$ws.Range("K58").Value = 7576235.399999999
$ws.Range("M58").Value = -7576085.399999999
$ws.Range("H86").Value = 6160.45
$ws.Range("I86").Value = 7241.75
$ws.Range("J86").Value = 5439.5835
$ws.Range("K86").Value = 7241.75
$ws.Range("L86").Value = 5439.5835
$ws.Range("M86").Value = -6118.75
$ws.Range("N86").Value = -7685.5835
$ws.Range("H87").Value = 37997
$ws.Range("I87").Value = 30800
$ws.Range("J87").Value = 39025.145
$ws.Range("K87").Value = 30800
$ws.Range("L87").Value = 39025.145
$ws.Range("M87").Value = -29552
$ws.Range("N87").Value = -41521.145
$ws.Range("H89").Value = 6160.45
$ws.Range("I89").Value = 7241.75
$ws.Range("J89").Value = 5439.5835
$ws.Range("K89").Value = 36208.75
$ws.Range("L89").Value = 27197.9175
$ws.Range("M89").Value = -30592.75
$ws.Range("N89").Value = -38429.9175
$ws.Range("H90").Value = 37997
$ws.Range("I90").Value = 30800
$ws.Range("J90").Value = 39025.145
$ws.Range("K90").Value = 92400
$ws.Range("L90").Value = 117075.435
$ws.Range("M90").Value = -86160
$ws.Range("N90").Value = -129555.435
$ws.Range("H116").Value = 1920.8125
$ws.Range("I116").Value = 1136.4615
$ws.Range("J116").Value = 5319.6665
$ws.Range("K116").Value = 1136.4615
$ws.Range("L116").Value = 5319.6665
$ws.Range("M116").Value = 2305.5385
$ws.Range("N116").Value = -12203.6665
$ws.Range("H138").Value = 4907.283
$ws.Range("J138").Value = 7016.625
$ws.Range("L138").Value = 21049.875
$ws.Range("N138").Value = -31329.875

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1745.4717
$ws.Range("I61").Value = 1122.3704
$ws.Range("J61").Value = 2392.5386
$ws.Range("K61").Value = 1122.3704
$ws.Range("L61").Value = 2392.5386
$ws.Range("M61").Value = -910.3704
$ws.Range("N61").Value = -2816.5386
$ws.Range("H63").Value = 2298.6365
$ws.Range("I63").Value = 1017
$ws.Range("K63").Value = 1017
$ws.Range("M63").Value = -331
$ws.Range("H66").Value = 2298.6365
$ws.Range("I66").Value = 1017
$ws.Range("K66").Value = 5085
$ws.Range("M66").Value = -1653
$ws.Range("H74").Value = 1837.28
$ws.Range("I74").Value = 1314.4
$ws.Range("J74").Value = 1968
$ws.Range("K74").Value = 1314.4
$ws.Range("L74").Value = 1968
$ws.Range("M74").Value = -440.4000000000001
$ws.Range("N74").Value = -3716
$ws.Range("H77").Value = 1837.28
$ws.Range("I77").Value = 1314.4
$ws.Range("J77").Value = 1968
$ws.Range("K77").Value = 6572
$ws.Range("L77").Value = 9840
$ws.Range("M77").Value = -2204
$ws.Range("N77").Value = -18576
$ws.Range("H132").Value = 2796.8538
$ws.Range("I132").Value = 2740.6667
$ws.Range("J132").Value = 3201.4
$ws.Range("K132").Value = 8222.000100000001
$ws.Range("L132").Value = 9604.200000000001
$ws.Range("M132").Value = -5692.000100000001
$ws.Range("N132").Value = -14664.2
$ws.Range("H136").Value = 1745.4717
$ws.Range("I136").Value = 1122.3704
$ws.Range("J136").Value = 2392.5386
$ws.Range("K136").Value = 3367.1112
$ws.Range("L136").Value = 7177.6158
$ws.Range("M136").Value = -817.1112000000003
$ws.Range("N136").Value = -12277.6158

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1763.3939
$ws.Range("I134").Value = 1818.6666
$ws.Range("K134").Value = 5455.9998
$ws.Range("M134").Value = -2920.9998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 9935.643
$ws.Range("J41").Value = 10695.385
$ws.Range("L41").Value = 10695.385
$ws.Range("N41").Value = -11551.385
$ws.Range("H58").Value = 1410.641
$ws.Range("I58").Value = 1168.4073
$ws.Range("J58").Value = 1955.6666
$ws.Range("K58").Value = 1168.4073
$ws.Range("L58").Value = 1955.6666
$ws.Range("M58").Value = -965.4073000000001
$ws.Range("N58").Value = -2361.6666
$ws.Range("H105").Value = 931.96155
$ws.Range("I105").Value = 883.2727
$ws.Range("J105").Value = 1199.75
$ws.Range("K105").Value = 883.2727
$ws.Range("L105").Value = 1199.75
$ws.Range("M105").Value = 863.7273
$ws.Range("N105").Value = -4693.75
$ws.Range("H107").Value = 718.26086
$ws.Range("I107").Value = 826.2727
$ws.Range("J107").Value = 619.25
$ws.Range("K107").Value = 826.2727
$ws.Range("L107").Value = 619.25
$ws.Range("M107").Value = 1093.7273
$ws.Range("N107").Value = -4459.25
$ws.Range("H136").Value = 1410.641
$ws.Range("I136").Value = 1168.4073
$ws.Range("J136").Value = 1955.6666
$ws.Range("K136").Value = 3505.2219
$ws.Range("L136").Value = 5866.9998
$ws.Range("M136").Value = -955.2219000000005
$ws.Range("N136").Value = -10966.9998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 900
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 900
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H62").Value = 4796.25
$ws.Range("J62").Value = 4796.25
$ws.Range("L62").Value = 14388.75
$ws.Range("N62").Value = -15760.75
$ws.Range("H65").Value = 4796.25
$ws.Range("J65").Value = 4796.25
$ws.Range("L65").Value = 43166.25
$ws.Range("N65").Value = -50030.25
$ws.Range("H68").Value = 19579.672
$ws.Range("I68").Value = 1092
$ws.Range("J68").Value = 25892.537
$ws.Range("K68").Value = 3276
$ws.Range("L68").Value = 77677.611
$ws.Range("M68").Value = -2465
$ws.Range("N68").Value = -79299.611
$ws.Range("H69").Value = 1569.9
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 1569.9
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 4709.700000000001
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -6331.700000000001
$ws.Range("H71").Value = 19579.672
$ws.Range("I71").Value = 1092
$ws.Range("J71").Value = 25892.537
$ws.Range("K71").Value = 9828
$ws.Range("L71").Value = 233032.833
$ws.Range("M71").Value = -5772
$ws.Range("N71").Value = -241144.833
$ws.Range("H72").Value = 1569.9
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 1569.9
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 14129.1
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -22241.1
$ws.Range("H103").Value = 813.7143
$ws.Range("I103").Value = 116
$ws.Range("J103").Value = 5000
$ws.Range("K103").Value = 348
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = 531
$ws.Range("N103").Value = -16758
$ws.Range("H107").Value = 1051.2444
$ws.Range("J107").Value = 1296.4642
$ws.Range("L107").Value = 3889.3926
$ws.Range("N107").Value = -7729.392599999999
$ws.Range("H131").Value = 2285.95
$ws.Range("I131").Value = 576.6667
$ws.Range("J131").Value = 2587.5881
$ws.Range("K131").Value = 1730.0001
$ws.Range("L131").Value = 7762.7643
$ws.Range("M131").Value = 3309.9999
$ws.Range("N131").Value = -17842.7643

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 9433.333000000001
$ws.Range("I58").Value = 1350
$ws.Range("K58").Value = 1350
$ws.Range("M58").Value = -1073
$ws.Range("H80").Value = 3750
$ws.Range("I80").Value = 5400
$ws.Range("J80").Value = 2100
$ws.Range("K80").Value = 5400
$ws.Range("L80").Value = 2100
$ws.Range("M80").Value = -4402
$ws.Range("N80").Value = -4096
$ws.Range("H83").Value = 3750
$ws.Range("I83").Value = 5400
$ws.Range("J83").Value = 2100
$ws.Range("K83").Value = 27000
$ws.Range("L83").Value = 10500
$ws.Range("M83").Value = -22008
$ws.Range("N83").Value = -20484
$ws.Range("H130").Value = 44980
$ws.Range("J130").Value = 44980
$ws.Range("L130").Value = 44980
$ws.Range("N130").Value = -55020
$ws.Range("H132").Value = 2651.2144
$ws.Range("I132").Value = 2406.75
$ws.Range("J132").Value = 3262.375
$ws.Range("K132").Value = 7220.25
$ws.Range("L132").Value = 9787.125
$ws.Range("M132").Value = -4690.25
$ws.Range("N132").Value = -14847.125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7422.1113
$ws.Range("I132").Value = 8216.666999999999
$ws.Range("J132").Value = 5833
$ws.Range("K132").Value = 24650.001
$ws.Range("L132").Value = 17499
$ws.Range("M132").Value = -22120.001
$ws.Range("N132").Value = -22559
$ws.Range("H136").Value = 1116.8422
$ws.Range("I136").Value = 956.2593000000001
$ws.Range("J136").Value = 1511
$ws.Range("K136").Value = 2868.7779
$ws.Range("L136").Value = 4533
$ws.Range("M136").Value = -318.7779
$ws.Range("N136").Value = -9633

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2698.3057
$ws.Range("I132").Value = 3691.5
$ws.Range("J132").Value = 1705.1111
$ws.Range("K132").Value = 11074.5
$ws.Range("L132").Value = 5115.3333
$ws.Range("M132").Value = -8544.5
$ws.Range("N132").Value = -10175.3333
$ws.Range("H136").Value = 815.5806
$ws.Range("I136").Value = 550.13635
$ws.Range("J136").Value = 1464.4445
$ws.Range("K136").Value = 1650.40905
$ws.Range("L136").Value = 4393.333500000001
$ws.Range("M136").Value = 899.59095
$ws.Range("N136").Value = -9493.333500000001
